$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, pushing existing rows 103-181 down to 104-182
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with its data
$ws.Range("A103").Value = 4
$ws.Range("B103").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C103").Value = "Los Lagos"
$ws.Range("D103").Value = 44574
$ws.Range("E103").Value = 10
$ws.Range("F103").Value = "Fruta"
$ws.Range("G103").Value = 100101
$ws.Range("H103").Value = "Berries"
$ws.Range("I103").Value = 100112025
$ws.Range("J103").Value = "Frutilla"
$ws.Range("K103").Value = "Sin especificar"
$ws.Range("L103").Value = "Primera"
$ws.Range("M103").Value = 400
$ws.Range("N103").Value = 8500
$ws.Range("O103").Value = 9000
$ws.Range("P103").Value = 8750
$ws.Range("Q103").Value = "`$/caja 7 kilos"
$ws.Range("R103").Value = "Región de La Araucanía"
$ws.Range("S103").Value = 1250
$ws.Range("T103").Value = 7
